$wb = $excel.ActiveWorkbook

# --- Sheet-3th: add header row "This is a test!", "aaa bbb ccc", "1234567890" ---
$ws3 = $wb.Worksheets.Item("Sheet-3th")

$ws3.Range("A1").Value = "This is a test!"
$ws3.Range("B1").Value = "aaa bbb ccc"

# C1 must hold "1234567890" as *text* (shared string), not a number, even
# though it's all digits. Forcing text via an apostrophe prefix (same as
# typing '1234567890 into a cell) stores it as a string, then clearing the
# formatting drops the transient quote-prefix style that operation adds so
# the cell keeps a plain string type with no explicit style index.
$ws3.Range("C1").Value = "'1234567890"
$ws3.Range("C1").ClearFormats()

# --- activate Sheet-Two, which flips workbookView activeTab to 1 and moves
# tabSelected from Sheet-One's sheetView to Sheet-Two's sheetView ---
$ws2 = $wb.Worksheets.Item("Sheet-Two")
$ws2.Activate()
